$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BB column: BB1 header date, then BB2:BB82 forecast values, BB83 ---

# BB1 = new date header (copy style/format from BA1, then set the new date value)
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# BB2:BB70 = copies of BA2:BA70 (same forecast values carried into new column)
$sameVals = @(
    1.63165117816655, 1.358419700277324, -0.1398451755242718, -0.07003068004286206,
    -7.3868590581191, -10.33903131837766, -1.060235412222937, 3.423433284297019,
    2.063081011733999, 2.99202665168275, 6.534207423404695, 1.294218807309846,
    2.007778863461724, 2.488825190489734, -0.02236378853282872, 1.864148130528193,
    0.3074416423962276, 1.591452302439862, 1.24229439238384, 1.021051258256691,
    -2.083516682089652, 0.8945075486844729, 1.373553924371535, 0.3276711086097635,
    1.809509533486136, 0.4339813219297071, 0.4041977135476031, 1.507492882068462,
    1.350099422002103, 1.114800352984872, 1.516137977472326, -0.02409266101658147,
    -0.7388860519741201, 1.567695497950282, 1.091382109524247, -0.2679292251141305,
    1.682996656587392, 1.591260107191601, 1.13706115148689, 1.735453665039003,
    2.625533283765208, -0.3243298885145123, 0.6372341859553217, -0.5789332341234967,
    0.1728728569232914, 1.7, -1.4, 1.3, -0.6, -3.292009884772611, -20.353166912592,
    17.87631681612835, 4.395171409529297, 1.763269038133103, 1.651360288740889,
    -0.6023499639201475, 3.773751128807561, 0.3486937127494798, 0.6021692060546258,
    1.505090528591751, -1.274647749701984, 0.4647476577308112, -0.6885297541760451,
    -0.9272934801906558, -0.8572418502401149, 0.05364604092510206, 1.830682919206694,
    -2.604576060482884, -3.124437332092583
)

# BB71:BB82 = newly computed forecast values (differ from BA71:BA82)
$newVals = @(
    2.411315004676197, 0.2886831937783967, -0.7196185376451893, 0.8547640426299213,
    0.5836696676700008, 0.3581087741383818, 0.4777118888332296, 0.486068368610706,
    0.4603442167010163, 0.4663515018921451, 0.4695689069109032, 0.4674468967903231
)

$allVals = $sameVals + $newVals   # rows 2..82 => 81 values total

$n = $allVals.Count
$data = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $data[$i,0] = $allVals[$i]
}
$ws.Range("BB2:BB82").Value = $data

# --- New row 83: A83 date (styled like column A header cells), BB83 value ---
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.4674511079923597
